$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

$cell = $t.Cell(1,1)
$cell.Range.Text = "57 x 87" + $br + "  8    7" + $br + "  ----" + $br + "5|    |" + $br + "7|    |"

$cell = $t.Cell(1,2)
$cell.Range.Text = "57 x 50" + $br + "  5    0" + $br + "  ----" + $br + "5|    |" + $br + "7|    |"

$cell = $t.Cell(1,3)
$cell.Range.Text = "18 x 29" + $br + "  2    9" + $br + "  ----" + $br + "1|    |" + $br + "8|    |"

$cell = $t.Cell(2,1)
$cell.Range.Text = "89 x 25" + $br + "  2    5" + $br + "  ----" + $br + "8|    |" + $br + "9|    |"

$cell = $t.Cell(2,2)
$cell.Range.Text = "73 x 57" + $br + "  5    7" + $br + "  ----" + $br + "7|    |" + $br + "3|    |"

$cell = $t.Cell(2,3)
$cell.Range.Text = "90 x 72" + $br + "  7    2" + $br + "  ----" + $br + "9|    |" + $br + "0|    |"

$cell = $t.Cell(3,1)
$cell.Range.Text = "94 x 52" + $br + "  5    2" + $br + "  ----" + $br + "9|    |" + $br + "4|    |"

$cell = $t.Cell(3,2)
$cell.Range.Text = "67 x 17" + $br + "  1    7" + $br + "  ----" + $br + "6|    |" + $br + "7|    |"

$cell = $t.Cell(3,3)
$cell.Range.Text = "34 x 45" + $br + "  4    5" + $br + "  ----" + $br + "3|    |" + $br + "4|    |"

$cell = $t.Cell(4,1)
$cell.Range.Text = "45 x 44" + $br + "  4    4" + $br + "  ----" + $br + "4|    |" + $br + "5|    |"

$cell = $t.Cell(4,2)
$cell.Range.Text = "15 x 37" + $br + "  3    7" + $br + "  ----" + $br + "1|    |" + $br + "5|    |"

$cell = $t.Cell(4,3)
$cell.Range.Text = "95 x 63" + $br + "  6    3" + $br + "  ----" + $br + "9|    |" + $br + "5|    |"

$cell = $t.Cell(5,1)
$cell.Range.Text = "16 x 98" + $br + "  9    8" + $br + "  ----" + $br + "1|    |" + $br + "6|    |"

$cell = $t.Cell(5,2)
$cell.Range.Text = "95 x 55" + $br + "  5    5" + $br + "  ----" + $br + "9|    |" + $br + "5|    |"

$cell = $t.Cell(5,3)
$cell.Range.Text = "91 x 59" + $br + "  5    9" + $br + "  ----" + $br + "9|    |" + $br + "1|    |"
